# "Signed Off time sheets"
#
# The supervisor (Ankita Gangotra) signs off the timesheet: her name and
# initials are filled in, the sign-off date is stamped, and the three
# days that were missing/short hours (Wed 5/14, Sat 5/17, Sun 5/18) are
# corrected to 3 hours each - the dependent totals/pay recalc
# automatically. Selection ends up resting on G26.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Supervisor name (merged G6:I6, mirrors the "Employee Name" field in G4) ---
$ws.Range("G6").Value = "Ankita Gangotra"
$ws.Range("G6").HorizontalAlignment = -4131   # xlLeft - reuse the existing named-field style

# --- Supervisor signature block (row 27): typed initials + sign-off date ---
$ws.Range("A27").Value = "A.G"
$ws.Range("A27").HorizontalAlignment = -4131  # xlLeft - reuse the existing named-field style

$ws.Range("D27").Value = 41800                 # 6/10/2014
$ws.Range("D27").NumberFormat = "mm-dd-yy"

# --- Corrected daily hours for Wed/Sat/Sun -----------------------------------
$ws.Range("B13").Value = 3
$ws.Range("B16").Value = 3
$ws.Range("B17").Value = 3

# --- Selection ends up on G26 (next to the "Date" label) --------------------
$ws.Range("G26").Select()
